# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 636.6667
$ws.Range("I28").Value = 636.3570999999999
$ws.Range("J28").Value = 637.75
$ws.Range("K28").Value = 636.3570999999999
$ws.Range("L28").Value = 637.75
$ws.Range("M28").Value = -151.3570999999999
$ws.Range("N28").Value = -1607.75
$ws.Range("H137").Value = 6237.7144
$ws.Range("I137").Value = 6733
$ws.Range("K137").Value = 20199
$ws.Range("M137").Value = -17649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 9666.666999999999
$ws.Range("J49").Value = 9666.666999999999
$ws.Range("L49").Value = 9666.666999999999
$ws.Range("N49").Value = -10186.667
$ws.Range("H74").Value = 307537.88
$ws.Range("I74").Value = 417747.97
$ws.Range("J74").Value = 87117.75
$ws.Range("K74").Value = 417747.97
$ws.Range("L74").Value = 87117.75
$ws.Range("M74").Value = -416873.97
$ws.Range("N74").Value = -88865.75
$ws.Range("H77").Value = 307537.88
$ws.Range("I77").Value = 417747.97
$ws.Range("J77").Value = 87117.75
$ws.Range("K77").Value = 2088739.85
$ws.Range("L77").Value = 435588.75
$ws.Range("M77").Value = -2084371.85
$ws.Range("N77").Value = -444324.75
$ws.Range("H80").Value = 21493
$ws.Range("J80").Value = 21493
$ws.Range("L80").Value = 21493
$ws.Range("N80").Value = -23489
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40722
$ws.Range("H83").Value = 21493
$ws.Range("J83").Value = 21493
$ws.Range("L83").Value = 64479
$ws.Range("N83").Value = -74463
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42496
$ws.Range("H86").Value = 39000
$ws.Range("J86").Value = 39000
$ws.Range("L86").Value = 39000
$ws.Range("N86").Value = -41372
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H89").Value = 39000
$ws.Range("J89").Value = 39000
$ws.Range("L89").Value = 117000
$ws.Range("N89").Value = -128856
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H132").Value = 30268.918
$ws.Range("I132").Value = 40255.406
$ws.Range("J132").Value = 3305.4
$ws.Range("K132").Value = 120766.218
$ws.Range("L132").Value = 9916.200000000001
$ws.Range("M132").Value = -118236.218
$ws.Range("N132").Value = -14976.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 78.5
$ws.Range("I80").Value = 103.28571
$ws.Range("J80").Value = 59.22222
$ws.Range("K80").Value = 103.28571
$ws.Range("L80").Value = 59.22222
$ws.Range("M80").Value = 894.71429
$ws.Range("N80").Value = -2055.22222
$ws.Range("H83").Value = 78.5
$ws.Range("I83").Value = 103.28571
$ws.Range("J83").Value = 59.22222
$ws.Range("K83").Value = 516.42855
$ws.Range("L83").Value = 296.1111
$ws.Range("M83").Value = 4475.57145
$ws.Range("N83").Value = -10280.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 873.5
$ws.Range("I16").Value = 771.0909
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 771.0909
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -484.0909
$ws.Range("N16").Value = -2574
$ws.Range("H31").Value = 3688
$ws.Range("I31").Value = 1890.3529
$ws.Range("K31").Value = 1890.3529
$ws.Range("M31").Value = -1595.3529
$ws.Range("H34").Value = 3688
$ws.Range("I34").Value = 1890.3529
$ws.Range("K34").Value = 1890.3529
$ws.Range("M34").Value = -1688.3529
$ws.Range("H35").Value = 111113160
$ws.Range("I35").Value = 142858350
$ws.Range("K35").Value = 142858350
$ws.Range("M35").Value = -142858056
$ws.Range("H113").Value = 873.5
$ws.Range("I113").Value = 771.0909
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 771.0909
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1398.9091
$ws.Range("N113").Value = -6340
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2566.5652
$ws.Range("I69").Value = 589.125
$ws.Range("K69").Value = 1767.375
$ws.Range("M69").Value = -956.375
$ws.Range("H72").Value = 2566.5652
$ws.Range("I72").Value = 589.125
$ws.Range("K72").Value = 5302.125
$ws.Range("M72").Value = -1246.125
$ws.Range("H131").Value = 2121.4707
$ws.Range("J131").Value = 1708.4067
$ws.Range("L131").Value = 5125.2201
$ws.Range("N131").Value = -15205.2201

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 599.3333
$ws.Range("I19").Value = 638
$ws.Range("J19").Value = 406
$ws.Range("K19").Value = 638
$ws.Range("L19").Value = 406
$ws.Range("M19").Value = -350
$ws.Range("N19").Value = -982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 1200
$ws.Range("J14").Value = 2200
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 2200
$ws.Range("M14").Value = -1028
$ws.Range("N14").Value = -2544
$ws.Range("H132").Value = 12507790
$ws.Range("I132").Value = 17860258
$ws.Range("J132").Value = 18699.5
$ws.Range("K132").Value = 53580774
$ws.Range("L132").Value = 56098.5
$ws.Range("M132").Value = -53578244
$ws.Range("N132").Value = -61158.5
$ws.Range("H136").Value = 4414.3784
$ws.Range("I136").Value = 2282.963
$ws.Range("J136").Value = 10169.2
$ws.Range("K136").Value = 6848.889000000001
$ws.Range("L136").Value = 30507.6
$ws.Range("M136").Value = -4298.889000000001
$ws.Range("N136").Value = -35607.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 21500
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17005
$ws.Range("H132").Value = 23257794
$ws.Range("I132").Value = 30304370
$ws.Range("J132").Value = 4088.8
$ws.Range("K132").Value = 90913110
$ws.Range("L132").Value = 12266.4
$ws.Range("M132").Value = -90910580
$ws.Range("N132").Value = -17326.4

Write-Host "Applied scheduled market-data refresh to all sheets."